$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that is bumped by one day
# (45171 -> 45172, i.e. 2023-09-02 -> 2023-09-03) for every data row (2..244).
$newDate = 45172

for ($row = 2; $row -le 244; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
